$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Level"
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 2

$ws.Range("F5").Select()
